$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextPlaceholder($targetRef, $donorRef) {
    $ws.Range($donorRef).Copy()
    $ws.Range($targetRef).PasteSpecial(-4122)
    $ws.Range($donorRef).Copy()
    $ws.Range($targetRef).PasteSpecial()
}

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 31   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  4/1/2024  Through  4/7/2024"

# --- Numeric cell updates ---
$ws.Range("L16").Value = 100
$ws.Range("N16").Value = -60
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 33
$ws.Range("J17").Value = 22
$ws.Range("K17").Value = 50
$ws.Range("L17").Value = 153.846153846154
$ws.Range("M17").Value = 175
$ws.Range("N17").Value = 83.333333333333
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = -25
$ws.Range("I18").Value = 6
$ws.Range("K18").Value = -68.421052631578
$ws.Range("L18").Value = 200
$ws.Range("M18").Value = -80
$ws.Range("N18").Value = -91.304347826087
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -44.444444444444
$ws.Range("F19").Value = 21
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = -16
$ws.Range("I19").Value = 72
$ws.Range("J19").Value = 75
$ws.Range("K19").Value = -4
$ws.Range("L19").Value = -4
$ws.Range("M19").Value = 67.441860465116
$ws.Range("N19").Value = 100
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 3
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 13
$ws.Range("J20").Value = 19
$ws.Range("K20").Value = -31.578947368421
$ws.Range("L20").Value = -53.571428571428
$ws.Range("M20").Value = 30
$ws.Range("N20").Value = -92.696629213483
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = -14.285714285714
$ws.Range("F21").Value = 36
$ws.Range("G21").Value = 44
$ws.Range("H21").Value = -18.181818181818
$ws.Range("I21").Value = 133
$ws.Range("J21").Value = 140
$ws.Range("K21").Value = -5
$ws.Range("L21").Value = 9.9173553719
$ws.Range("M21").Value = 35.714285714285
$ws.Range("N21").Value = -58.307210031348
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = -44.444444444444
$ws.Range("F24").Value = 22
$ws.Range("G24").Value = 42
$ws.Range("H24").Value = -47.619047619047
$ws.Range("I24").Value = 103
$ws.Range("J24").Value = 144
$ws.Range("K24").Value = -28.472222222222
$ws.Range("L24").Value = 1.980198019801
$ws.Range("M24").Value = -12.711864406779
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = -19.047619047619
$ws.Range("I25").Value = 53
$ws.Range("J25").Value = 83
$ws.Range("K25").Value = -36.144578313253
$ws.Range("L25").Value = 140.909090909091
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 11
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = -31.25
$ws.Range("I26").Value = 36
$ws.Range("J26").Value = 54
$ws.Range("L26").Value = -29.411764705882
$ws.Range("M26").Value = -30.76923076923
$ws.Range("L28").Value = -75

# --- Text-placeholder (str20/str21) cell updates ---
# donors: D14 has style14+"0"(str20); E14 has style14+"***.*"(str21)
Set-TextPlaceholder "G14" "D14"
Set-TextPlaceholder "H14" "E14"
Set-TextPlaceholder "D18" "D14"
Set-TextPlaceholder "E18" "E14"
Set-TextPlaceholder "C28" "D14"
Set-TextPlaceholder "C31" "D14"
